$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "Nazarena Raos"
$ws.Range("B31").Value = "Alberto Cerisara | SHARK ATTACK"
$ws.Range("C31").Value = "Andrea Gober | U.SGUARNA"
$ws.Range("D31").Value = "Roberto Barozzi | demobusters"
$ws.Range("E31").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("F31").Value = "Gabriele Lasta | RSA United"
